$d = $word.ActiveDocument

$replacements = @(
    @("2024-07-20 Saturday", "2024-07-21 Sunday"),
    @("128×2=", "145×6="),
    @("577×2=", "151×9="),
    @("502×8=", "427×4="),
    @("765×6=", "777×2="),
    @("327×5=", "317×7="),
    @("639×4=", "925×4="),
    @("127×3=", "900×7="),
    @("893×8=", "483×3="),
    @("108×8=", "621×4="),
    @("760×5=", "881×7="),
    @("381×6=", "208×5="),
    @("772×2=", "922×4="),
    @("293×6=", "376×8="),
    @("793×9=", "208×3="),
    @("407×2=", "976×6="),
    @("168×8=", "403×8="),
    @("983×8=", "312×8="),
    @("437×5=", "995×3="),
    @("423×7=", "955×6="),
    @("272×3=", "250×6="),
    @("430×4=", "984×4="),
    @("825×2=", "748×4="),
    @("669×4=", "198×7="),
    @("765×8=", "530×6="),
    @("773×2=", "340×3=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
